$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------

# Row 7 "stage": stored text "2" becomes a real number 2
$ws.Range("B7").Value = 2

# Row 8 "dayUpdate": refresh the timestamp
$ws.Range("B8").Value = "07/28/2023 12:20:09"

# Row 9 "uploadFolderName": refresh the folder name
$ws.Range("B9").Value = "28-07-2023"

# New row 12: mailPhuTrach
$ws.Range("A12").Value = "mailPhuTrach"
$ws.Range("A12").VerticalAlignment = -4108   # xlCenter

$ws.Range("B12").Value = "kemclone1@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:kemclone1@gmail.com") | Out-Null
$ws.Range("B12").WrapText = $true
$ws.Range("B12").VerticalAlignment = -4108   # xlCenter

# New row 13: mailKhachHang
$ws.Range("A13").Value = "mailKhachHang"
$ws.Range("A13").VerticalAlignment = -4108   # xlCenter

$ws.Range("B13").Value = "kemclone3@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:kemclone3@gmail.com") | Out-Null
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4108   # xlCenter

# --- Formatting updates ------------------------------------------------

# Base font + hyperlink font: Calibri -> Arial
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Arial"
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Font.Name = "Arial"

# Restore the header row's theme fill (changing the base font above can
# reset derived cell formatting), keep the same look as before the font swap
$ws.Range("A1:B1").Interior.Pattern = 1
$ws.Range("A1:B1").Interior.ThemeColor = 5
$ws.Range("A1:B1").Interior.TintAndShade = 0

# --- Selection cosmetics ------------------------------------------------
$ws.Range("B7").Select() | Out-Null

Write-Output "done"
